$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in actual clock in/out times that were recorded for this pay period
$ws.Range("C16").Value = 0.58333333333333337
$ws.Range("F16").Value = 0.64583333333333337

$ws.Range("C19").Value = 0.5
$ws.Range("F19").Value = 0.66666666666666663

# Grand-total cell no longer needed here
$ws.Range("H21").ClearContents()

# Scroll the view down and move the selection
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H22").Select()
